$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet (Last Updated timestamp) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 11:10 AM"

# --- Update "distance from Dma50" sheet column C values ---
$ws = $wb.Worksheets.Item("distance from Dma50")

$ws.Range("C2").Value = 9.9132
$ws.Range("C3").Value = 7.3672
$ws.Range("C4").Value = 6.349
$ws.Range("C5").Value = 5.1845
$ws.Range("C6").Value = 5.1432
$ws.Range("C7").Value = 5.0608
$ws.Range("C8").Value = 4.4862
$ws.Range("C9").Value = 4.3707
$ws.Range("C10").Value = 3.8378
$ws.Range("C11").Value = 3.5222
$ws.Range("C12").Value = 3.5081
$ws.Range("C13").Value = 3.3496
$ws.Range("C14").Value = 3.1423
$ws.Range("C15").Value = 3.0807
$ws.Range("C16").Value = 3.015
$ws.Range("C17").Value = 2.853
$ws.Range("C18").Value = 2.7119
$ws.Range("C19").Value = 2.6529
$ws.Range("C20").Value = 2.3495
$ws.Range("C21").Value = 2.3041
$ws.Range("C22").Value = 1.432
$ws.Range("C23").Value = 1.4017
$ws.Range("C24").Value = 1.3976
$ws.Range("C25").Value = 1.0649
$ws.Range("C26").Value = 1.0431
$ws.Range("C27").Value = 0.856
$ws.Range("C28").Value = 0.8120000000000001
$ws.Range("C29").Value = 0.2507
$ws.Range("C30").Value = -2.1113
